$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2128712871287129
$ws.Range("C2").Value = 0.5247524752475248
$ws.Range("J2").Value = 0.01485148514851485
$ws.Range("P2").Value = 0.1287128712871287
$ws.Range("S2").Value = 0.1188118811881188
$ws.Range("B3").Value = 0.01851851851851852
$ws.Range("C3").Value = 0.02777777777777778
$ws.Range("J3").Value = 0.04629629629629629
$ws.Range("P3").Value = 0.7777777777777778
$ws.Range("S3").Value = 0.1296296296296296
$ws.Range("J4").Value = 0.03571428571428571
$ws.Range("P4").Value = 0.6428571428571429
$ws.Range("S4").Value = 0.3214285714285715
$ws.Range("B6").Value = 0.09659090909090909
$ws.Range("D6").Value = 0.005681818181818182
$ws.Range("F6").Value = 0.05681818181818182
$ws.Range("J6").Value = 0.2443181818181818
$ws.Range("O6").Value = 0.01704545454545454
$ws.Range("Q6").Value = 0.1590909090909091
$ws.Range("R6").Value = 0.0625
$ws.Range("S6").Value = 0.3579545454545455
$ws.Range("B7").Value = 0.09219858156028368
$ws.Range("D7").Value = 0.02127659574468085
$ws.Range("F7").Value = 0.05673758865248227
$ws.Range("J7").Value = 0.1276595744680851
$ws.Range("O7").Value = 0.007092198581560284
$ws.Range("Q7").Value = 0.2127659574468085
$ws.Range("R7").Value = 0.09929078014184398
$ws.Range("S7").Value = 0.3829787234042553
$ws.Range("B8").Value = 0.08977556109725686
$ws.Range("D8").Value = 0.01246882793017456
$ws.Range("E8").Value = 0.002493765586034913
$ws.Range("F8").Value = 0.04488778054862843
$ws.Range("J8").Value = 0.1072319201995012
$ws.Range("O8").Value = 0.007481296758104738
$ws.Range("Q8").Value = 0.1920199501246883
$ws.Range("R8").Value = 0.09476309226932668
$ws.Range("S8").Value = 0.4488778054862843
$ws.Range("B9").Value = 0.06626506024096386
$ws.Range("E9").Value = 0.006024096385542169
$ws.Range("F9").Value = 0.07228915662650602
$ws.Range("J9").Value = 0.08433734939759036
$ws.Range("O9").Value = 0.01204819277108434
$ws.Range("Q9").Value = 0.1686746987951807
$ws.Range("R9").Value = 0.1385542168674699
$ws.Range("S9").Value = 0.4518072289156627
$ws.Range("B10").Value = 0.07766990291262135
$ws.Range("D10").Value = 0.01844660194174757
$ws.Range("E10").Value = 0.0009708737864077669
$ws.Range("F10").Value = 0.07087378640776699
$ws.Range("J10").Value = 0.116504854368932
$ws.Range("O10").Value = 0.007766990291262136
$ws.Range("Q10").Value = 0.2184466019417476
$ws.Range("R10").Value = 0.1145631067961165
$ws.Range("S10").Value = 0.374757281553398
$ws.Range("G11").Value = 0.1818181818181818
$ws.Range("J11").Value = 0.05194805194805195
$ws.Range("K11").Value = 0.2164502164502164
$ws.Range("L11").Value = 0.5324675324675324
$ws.Range("S11").Value = 0.01731601731601732
$ws.Range("G12").Value = 0.6904761904761905
$ws.Range("J12").Value = 0.253968253968254
$ws.Range("S12").Value = 0.05555555555555555
$ws.Range("G13").Value = 0.5806451612903226
$ws.Range("J13").Value = 0.3870967741935484
$ws.Range("S13").Value = 0.03225806451612903
$ws.Range("F15").Value = 0.03448275862068965
$ws.Range("H15").Value = 0.1724137931034483
$ws.Range("I15").Value = 0.08045977011494253
$ws.Range("J15").Value = 0.4022988505747127
$ws.Range("K15").Value = 0.06321839080459771
$ws.Range("M15").Value = 0.02298850574712644
$ws.Range("O15").Value = 0.08045977011494253
$ws.Range("S15").Value = 0.1436781609195402
$ws.Range("F16").Value = 0.01626016260162602
$ws.Range("H16").Value = 0.2276422764227642
$ws.Range("I16").Value = 0.1056910569105691
$ws.Range("J16").Value = 0.3902439024390244
$ws.Range("K16").Value = 0.1138211382113821
$ws.Range("M16").Value = 0.01626016260162602
$ws.Range("O16").Value = 0.03252032520325204
$ws.Range("S16").Value = 0.0975609756097561
$ws.Range("F17").Value = 0.01030927835051546
$ws.Range("H17").Value = 0.1855670103092784
$ws.Range("I17").Value = 0.07474226804123711
$ws.Range("J17").Value = 0.4407216494845361
$ws.Range("K17").Value = 0.09536082474226804
$ws.Range("M17").Value = 0.02577319587628866
$ws.Range("O17").Value = 0.05927835051546392
$ws.Range("S17").Value = 0.1082474226804124
$ws.Range("F18").Value = 0.0198019801980198
$ws.Range("H18").Value = 0.202970297029703
$ws.Range("I18").Value = 0.1138613861386139
$ws.Range("J18").Value = 0.3564356435643564
$ws.Range("K18").Value = 0.0891089108910891
$ws.Range("M18").Value = 0.009900990099009901
$ws.Range("O18").Value = 0.09405940594059406
$ws.Range("S18").Value = 0.1138613861386139
$ws.Range("F19").Value = 0.01360544217687075
$ws.Range("H19").Value = 0.2264334305150632
$ws.Range("I19").Value = 0.08454810495626822
$ws.Range("J19").Value = 0.3673469387755102
$ws.Range("K19").Value = 0.09815354713313897
$ws.Range("M19").Value = 0.01360544217687075
$ws.Range("O19").Value = 0.07677356656948493
$ws.Range("S19").Value = 0.119533527696793
